$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.666.91"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.602.73"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.47"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.25"
$ws.Range("E8").Value = "  +8.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.832.29"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "1.603.11"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.555"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "29.688.96"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.14"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.17"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.09"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").Value = "0.0₃0698"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.49"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.71"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.48"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0480"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "1.422.80"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.86"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0496"
$ws.Range("E42").Value = "  +5.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.820"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  +20.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.47"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.39"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").Value = "1.741.56"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.69"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0526"
$ws.Range("E51").Value = "  +0.88%  "
